$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Overwrite rows 2-18 with the recalculated forecast series
# (old row N+1 data shifted up into row N, with recomputed y_1_forecast values)

$ws.Range("A2").Value = 39765
$ws.Range("B2").Value = 2008
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 2009
$ws.Range("E2").Value = 2.741420124134053

$ws.Range("A3").Value = 40130
$ws.Range("B3").Value = 2009
$ws.Range("C3").Value = -3.872359107260159
$ws.Range("D3").Value = 2010
$ws.Range("E3").Value = -1.796159642965267

$ws.Range("A4").Value = 40494
$ws.Range("B4").Value = 2010
$ws.Range("C4").Value = 4.530477057343663
$ws.Range("D4").Value = 2011
$ws.Range("E4").Value = 3.372216600811506

$ws.Range("A5").Value = 40862
$ws.Range("B5").Value = 2011
$ws.Range("C5").Value = 6.833902841285977
$ws.Range("D5").Value = 2012
$ws.Range("E5").Value = 5.794553957309168

$ws.Range("A6").Value = 41228
$ws.Range("B6").Value = 2012
$ws.Range("C6").Value = 4.166536506645224
$ws.Range("D6").Value = 2013
$ws.Range("E6").Value = 4.453225076234824

$ws.Range("A7").Value = 41592
$ws.Range("B7").Value = 2013
$ws.Range("C7").Value = 2.669880057548091
$ws.Range("D7").Value = 2014
$ws.Range("E7").Value = 4.443665867214519

$ws.Range("A8").Value = 41957
$ws.Range("B8").Value = 2014
$ws.Range("C8").Value = 5.50293301232252
$ws.Range("D8").Value = 2015
$ws.Range("E8").Value = 4.93942156312468

$ws.Range("A9").Value = 42321
$ws.Range("B9").Value = 2015
$ws.Range("C9").Value = 4.829481320500406
$ws.Range("D9").Value = 2016
$ws.Range("E9").Value = 5.304093839766777

$ws.Range("A10").Value = 42689
$ws.Range("B10").Value = 2016
$ws.Range("C10").Value = 5.100281927437122
$ws.Range("D10").Value = 2017
$ws.Range("E10").Value = 4.773858406282372

$ws.Range("A11").Value = 43053
$ws.Range("B11").Value = 2017
$ws.Range("C11").Value = 5.161358932333737
$ws.Range("D11").Value = 2018
$ws.Range("E11").Value = 5.448823958443616

$ws.Range("A12").Value = 43418
$ws.Range("B12").Value = 2018
$ws.Range("C12").Value = 5.902681694119694
$ws.Range("D12").Value = 2019
$ws.Range("E12").Value = 5.03659417252571

$ws.Range("A13").Value = 43783
$ws.Range("B13").Value = 2019
$ws.Range("C13").Value = 3.884502719230132
$ws.Range("D13").Value = 2020
$ws.Range("E13").Value = 4.009670676786059

$ws.Range("A14").Value = 44159
$ws.Range("B14").Value = 2020
$ws.Range("C14").Value = -3.840397826549158
$ws.Range("D14").Value = 2021
$ws.Range("E14").Value = -0.236010050592228

$ws.Range("A15").Value = 44525
$ws.Range("B15").Value = 2021
$ws.Range("C15").Value = 0.4839811651348835
$ws.Range("D15").Value = 2022
$ws.Range("E15").Value = 2.177145583294293

$ws.Range("A16").Value = 44890
$ws.Range("B16").Value = 2022
$ws.Range("C16").Value = 2.06342951900429
$ws.Range("D16").Value = 2023
$ws.Range("E16").Value = 0.3896432785800652

$ws.Range("A17").Value = 45254
$ws.Range("B17").Value = 2023
$ws.Range("C17").Value = -2.156362896191677
$ws.Range("D17").Value = 2024
$ws.Range("E17").Value = -1.437335768580206

$ws.Range("A18").Value = 45618
$ws.Range("B18").Value = 2024
$ws.Range("C18").Value = -0.8205034771073372
$ws.Range("D18").Value = 2025
$ws.Range("E18").Value = 0.07297157746815053

# Remove the now-obsolete last row (row 19), shrinking the used range to A1:E18
$ws.Rows("19").Delete()
